# The "Id" rows in this worksheet (rows 2-10) were reshuffled: each row
# is replaced wholesale by the data that used to live in a different row
# (row 11 / header row 1 are untouched). Re-apply every cell that differs
# between the original row and its new content, cell by cell, so that
# text-looking values (counts, dates, times) stay text instead of being
# auto-converted by Excel, and cells that become blank are cleared.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 now holds the data that was previously in row 10
$ws.Range("A2").Value = 111486400
$ws.Range("I2").Value = '''5'
$ws.Range("K2").Value = 'fullt utvecklade blad'
$ws.Range("Q2").Value = 624030.1824148977
$ws.Range("R2").Value = 6932961.620511409
$ws.Range("Y2").Value = '''2023-08-14'
$ws.Range("Z2").Value = '''00:00'
$ws.Range("AA2").Value = '''2023-08-14'
$ws.Range("AB2").Value = '''00:00'
$ws.Range("AC2").ClearContents()

# Row 3 now holds the data that was previously in row 2
$ws.Range("A3").Value = 111486347
$ws.Range("B3").Value = 96348
$ws.Range("D3").Value = 'VU'
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = 'Knärot'
$ws.Range("G3").Value = 'Goodyera repens'
$ws.Range("H3").Value = '(L.) R. Br.'
$ws.Range("I3").Value = '''20'
$ws.Range("J3").Value = 'plantor/tuvor'
$ws.Range("K3").Value = 'i frukt'
$ws.Range("N3").Value = 'observerad'
$ws.Range("P3").Value = 'Glödenhöjden (Glödenhöjden), Mpd'
$ws.Range("Q3").Value = 624030.5648888731
$ws.Range("R3").Value = 6933013.425735661
$ws.Range("S3").Value = 10
$ws.Range("Z3").Value = '''14:17'
$ws.Range("AB3").Value = '''14:17'
$ws.Range("AC3").Value = 'Ca 20 ex, flesta som bladrosetter. 1 överblommad fruktbildande'
$ws.Range("AH3").ClearContents()

# Row 4 now holds the data that was previously in row 5
$ws.Range("A4").Value = 111485917
$ws.Range("I4").Value = '''3'
$ws.Range("P4").Value = 'Glödenhöjden nordost (Glödenhöjden), Mpd'
$ws.Range("Q4").Value = 624090.1097011974
$ws.Range("R4").Value = 6933043.392863069
$ws.Range("S4").Value = 25
$ws.Range("Z4").Value = '''15:00'
$ws.Range("AB4").Value = '''15:00'
$ws.Range("AC4").ClearContents()

# Row 5 now holds the data that was previously in row 3
$ws.Range("A5").Value = 111485854
$ws.Range("B5").Value = 89405
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 1202
$ws.Range("F5").Value = 'Ullticka'
$ws.Range("G5").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H5").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I5").Value = '''1'
$ws.Range("J5").Value = 'fruktkroppar'
$ws.Range("K5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("Q5").Value = 624096.1730324102
$ws.Range("R5").Value = 6933042.231978768
$ws.Range("S5").Value = 20
$ws.Range("Z5").Value = '''15:02'
$ws.Range("AB5").Value = '''15:02'
$ws.Range("AH5").Value = 'Häll- och blockskog'

# Row 6 now holds the data that was previously in row 8
$ws.Range("A6").Value = 111486280
$ws.Range("I6").Value = '''3'
$ws.Range("J6").ClearContents()
$ws.Range("Q6").Value = 624009.7035872869
$ws.Range("R6").Value = 6933014.034667149
$ws.Range("Z6").Value = '''14:17'
$ws.Range("AB6").Value = '''14:17'
$ws.Range("AC6").ClearContents()

# Row 7 now holds the data that was previously in row 9
$ws.Range("A7").Value = 111486415
$ws.Range("I7").Value = '''10'
$ws.Range("J7").Value = 'plantor/tuvor'
$ws.Range("Q7").Value = 624040.2038791699
$ws.Range("R7").Value = 6932953.67081845
$ws.Range("Z7").Value = '''13:46'
$ws.Range("AB7").Value = '''13:46'

# Row 8 now holds the data that was previously in row 7
$ws.Range("A8").Value = 111486450
$ws.Range("I8").Value = '''2'
$ws.Range("K8").Value = 'fullt utvecklade blad'
$ws.Range("Q8").Value = 624051.1502826829
$ws.Range("R8").Value = 6932945.755648845
$ws.Range("Z8").Value = '''13:43'
$ws.Range("AB8").Value = '''13:43'

# Row 9 now holds the data that was previously in row 4
$ws.Range("A9").Value = 111486117
$ws.Range("K9").Value = 'blomning'
$ws.Range("Q9").Value = 623993.6707231236
$ws.Range("R9").Value = 6933021.760048959
$ws.Range("S9").Value = 15
$ws.Range("Z9").Value = '''14:17'
$ws.Range("AB9").Value = '''14:17'
$ws.Range("AC9").Value = '10 plantor varav 2 blommande'

# Row 10 now holds the data that was previously in row 6
$ws.Range("A10").Value = 111486385
$ws.Range("I10").Value = '''20'
$ws.Range("K10").Value = 'blomning'
$ws.Range("Q10").Value = 624029.7289273632
$ws.Range("R10").Value = 6932998.597210908
$ws.Range("Y10").Value = '''2023-08-12'
$ws.Range("Z10").Value = '''14:11'
$ws.Range("AA10").Value = '''2023-08-12'
$ws.Range("AB10").Value = '''14:11'
$ws.Range("AC10").Value = 'Ca 20 ex varav 3 blommande'
